# The data rows in this "invoice items" sheet were off-by-one: column C
# (the "tag" column, e.g. "D") was blank and the real tag/value/invoice#/
# subtotal/item-name data had spilled one column too far right, into an
# extra column G that the table header doesn't define.
#
# Fix: for every affected item row, shift the D:G block one column left
# into C:F (value AND formatting), then drop the now-unused column G cell
# entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(4, 5, 9, 13, 17, 21, 25, 29, 30, 31, 32, 36, 37, 38, 42, 43, 47, 51)

foreach ($r in $rows) {
    $src = $ws.Range("D" + $r + ":G" + $r)
    $dst = $ws.Range("C" + $r)
    $src.Copy($dst)
    $ws.Cells.Item($r, 7).Clear()
}
